$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Almacen" -> "Warehouse" (also pick up an explicit/automatic font color,
# distinguishing this cell's style from the sheet default)
$ws.Range("C2").Value = "Warehouse"
$ws.Range("C2").Font.Color = 0

# "Piso Productivo" -> "Fastenal"
$ws.Range("C3").Value = "Fastenal"

# Restore the cursor/selection to where the author left it on save
$ws.Range("G7").Select() | Out-Null
